# This script applies the "gh-pages output generated at 456a3b4" refresh:
# updated "want to go" / "min price" counters across all four sheets, plus
# the removal of two cancelled/duplicate shows from the "演出" (Shows) sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions): refresh "想去人数" counts ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 997
$wsExpo.Range("F5").Value = 446
$wsExpo.Range("F6").Value = 706
$wsExpo.Range("F9").Value = 27
$wsExpo.Range("F10").Value = 393
$wsExpo.Range("F11").Value = 201
$wsExpo.Range("F12").Value = 72
$wsExpo.Range("F13").Value = 809
$wsExpo.Range("F15").Value = 1973
$wsExpo.Range("F16").Value = 469
$wsExpo.Range("F17").Value = 6946
$wsExpo.Range("F18").Value = 511
$wsExpo.Range("F19").Value = 518
$wsExpo.Range("F21").Value = 88
$wsExpo.Range("F23").Value = 209

# ---- Sheet "演出" (Shows): drop the two 09-01 "音阅派国漫演唱会" rows ----
# (one was a duplicate already covered by another row, the other was cancelled);
# every later show shifts up by two rows and the remaining rows get refreshed
# "想去人数"/"最低票价" numbers, matching the freshly scraped bilibili data.
$wsShows = $wb.Worksheets.Item("演出")
$wsShows.Rows.Item(4).Delete()
$wsShows.Rows.Item(4).Delete()

# Renumber the leading index column (0-based row counter) after the delete
for ($i = 1; $i -le 18; $i++) {
    $wsShows.Cells.Item($i, 1).Value = $i - 1
}

# One more refreshed counter beyond the plain row-shift
$wsShows.Range("F11").Value = 53

# ---- Sheet "本地生活" (Local life): refresh counters ----
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5458
$wsLocal.Range("G2").Value = 30
$wsLocal.Range("F3").Value = 386
$wsLocal.Range("F4").Value = 380

# ---- Sheet "全部类型" (All types): refresh counters ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5458
$wsAll.Range("G3").Value = 30
$wsAll.Range("F4").Value = 386
$wsAll.Range("F5").Value = 380
$wsAll.Range("F8").Value = 31
$wsAll.Range("G8").Value = "不可售"
$wsAll.Range("F11").Value = 997
$wsAll.Range("F15").Value = 446
$wsAll.Range("F16").Value = 706
$wsAll.Range("F20").Value = 27
$wsAll.Range("F21").Value = 393
$wsAll.Range("F22").Value = 201
$wsAll.Range("F24").Value = 72
$wsAll.Range("F26").Value = 809
$wsAll.Range("F29").Value = 1973
$wsAll.Range("F30").Value = 469
$wsAll.Range("F31").Value = 6946
$wsAll.Range("F32").Value = 53
$wsAll.Range("F33").Value = 512
$wsAll.Range("F34").Value = 518
$wsAll.Range("F36").Value = 88
$wsAll.Range("F39").Value = 209
